# ---------------------------------------------------------------------------
# 1. Update the time_taken timestamps on the "data" sheet (column F, rows 2-25)
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:19:17.267059",
    "2021-10-05 14:19:17.267067",
    "2021-10-05 14:19:17.267070",
    "2021-10-05 14:19:17.267073",
    "2021-10-05 14:19:17.267076",
    "2021-10-05 14:19:17.267078",
    "2021-10-05 14:19:17.267081",
    "2021-10-05 14:19:17.267083",
    "2021-10-05 14:19:17.267086",
    "2021-10-05 14:19:17.267089",
    "2021-10-05 14:19:17.267091",
    "2021-10-05 14:19:17.267094",
    "2021-10-05 14:19:17.267096",
    "2021-10-05 14:19:17.267099",
    "2021-10-05 14:19:17.267101",
    "2021-10-05 14:19:17.267104",
    "2021-10-05 14:19:17.267107",
    "2021-10-05 14:19:17.267109",
    "2021-10-05 14:19:17.267112",
    "2021-10-05 14:19:17.267115",
    "2021-10-05 14:19:17.267117",
    "2021-10-05 14:19:17.267120",
    "2021-10-05 14:19:17.267122",
    "2021-10-05 14:19:17.267125"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" worksheet right after the "data" sheet
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1) - styled like the header row on the "data" sheet
# (bold font, thin border, centered horizontally, top vertically)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Index cell A2 - styled the same way as column A on the "data" sheet
$a2 = $metaSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

# Data row 2
$metaSheet.Range("B2").Value = "Bardet Biedl syndrome"
$metaSheet.Range("C2").Value = 543
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.11"
$metaSheet.Range("E2").Value = "2021-04-07T17:01:54.908533Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:17.263373"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/543/?format=json"

# Keep the "data" sheet as the active/selected tab, as before the edit
$dataSheet.Activate()
